$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.3
$ws.Range("E2").Value = 1.19
$ws.Range("D3").Value = 1.4
$ws.Range("F4").Value = 1.12
$ws.Range("B5").Value = 1.53
$ws.Range("C5").Value = 1.35
$ws.Range("D6").Value = 1.49
$ws.Range("E6").Value = 1.33
$ws.Range("G6").Value = 0.97
$ws.Range("F7").Value = 1.49
